# Clear the values in F5:H9 (for each egg definition row), leaving the
# cell styles intact. These cells previously held sample weight values
# (1, 2, 3) that are no longer wanted.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gacha")

$ws.Range("F5:H9").ClearContents()
